# Rename columns in the "vehicleTypes" and "shipments" tables of the
# transport sample-data workbook.
#
#   vehicleTypes:
#     maxCapacityLoadingMeter -> maxCapacityPallets
#     costsPerKm              -> costsPerDistanceUnit
#
#   shipments:
#     sender               -> senderId
#     senderServiceTime    -> senderStopDuration
#     recipient            -> recipientId
#     recipientServiceTime -> recipientStopDuration
#     loadingMeter         -> pallets
#     opportunityCosts     -> externalCosts

$wb = $excel.ActiveWorkbook

# --- vehicleTypes sheet -----------------------------------------------
$wsVehicleTypes = $wb.Worksheets.Item("vehicleTypes")
$wsVehicleTypes.Range("L1").Value = "maxCapacityPallets"
$wsVehicleTypes.Range("O1").Value = "costsPerDistanceUnit"

# --- shipments sheet ----------------------------------------------------
$wsShipments = $wb.Worksheets.Item("shipments")
$wsShipments.Range("C1").Value = "senderId"
$wsShipments.Range("G1").Value = "recipientId"
$wsShipments.Range("H1").Value = "recipientStopDuration"
$wsShipments.Range("M1").Value = "pallets"
$wsShipments.Range("N1").Value = "externalCosts"
$wsShipments.Range("D1").Value = "senderStopDuration"

# --- restore the active selection on each sheet -------------------------
$wsVehicleTypes.Activate() | Out-Null
$wsVehicleTypes.Range("O1").Select() | Out-Null

$wsShipments.Activate() | Out-Null
$wsShipments.Range("N1").Select() | Out-Null
